$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3423.9023
$ws.Range("I64").Value = 3046
$ws.Range("K64").Value = 3046
$ws.Range("M64").Value = -2798
$ws.Range("H67").Value = 3423.9023
$ws.Range("I67").Value = 3046
$ws.Range("K67").Value = 3046
$ws.Range("M67").Value = -2188
$ws.Range("H74").Value = 3698.889
$ws.Range("J74").Value = 3798.8235
$ws.Range("L74").Value = 3798.8235
$ws.Range("N74").Value = -5670.8235
$ws.Range("H76").Value = 3569.18
$ws.Range("I76").Value = 2959.2104
$ws.Range("K76").Value = 2959.2104
$ws.Range("M76").Value = -2644.2104
$ws.Range("H77").Value = 3698.889
$ws.Range("J77").Value = 3798.8235
$ws.Range("L77").Value = 18994.1175
$ws.Range("N77").Value = -28354.1175
$ws.Range("H79").Value = 3569.18
$ws.Range("I79").Value = 2959.2104
$ws.Range("K79").Value = 2959.2104
$ws.Range("M79").Value = -1867.2104
$ws.Range("H129").Value = 866.2174
$ws.Range("J129").Value = 875.0222
$ws.Range("L129").Value = 2625.0666
$ws.Range("N129").Value = -12625.0666
$ws.Range("H138").Value = 2381.08
$ws.Range("I138").Value = 1374.1555
$ws.Range("J138").Value = 3204.9272
$ws.Range("K138").Value = 4122.4665
$ws.Range("L138").Value = 9614.7816
$ws.Range("M138").Value = 1017.5335
$ws.Range("N138").Value = -19894.7816

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 8433.333000000001
$ws.Range("I31").Value = 8433.333000000001
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 8433.333000000001
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -8139.333000000001
$ws.Range("N31").Value = ""
$ws.Range("H32").Value = 5516.346
$ws.Range("I32").Value = 4718.8354
$ws.Range("K32").Value = 4718.8354
$ws.Range("M32").Value = -4431.8354
$ws.Range("H61").Value = 2642.3684
$ws.Range("I61").Value = 2402.2727
$ws.Range("J61").Value = 2972.5
$ws.Range("K61").Value = 2402.2727
$ws.Range("L61").Value = 2972.5
$ws.Range("M61").Value = -2190.2727
$ws.Range("N61").Value = -3396.5
$ws.Range("H122").Value = 2292.2727
$ws.Range("J122").Value = 2985
$ws.Range("L122").Value = 8955
$ws.Range("N122").Value = -13855
$ws.Range("H132").Value = 1375.541
$ws.Range("I132").Value = 871.35
$ws.Range("J132").Value = 2335.9048
$ws.Range("K132").Value = 2614.05
$ws.Range("L132").Value = 7007.714399999999
$ws.Range("M132").Value = -84.05000000000018
$ws.Range("N132").Value = -12067.7144
$ws.Range("H136").Value = 2642.3684
$ws.Range("I136").Value = 2402.2727
$ws.Range("J136").Value = 2972.5
$ws.Range("K136").Value = 7206.8181
$ws.Range("L136").Value = 8917.5
$ws.Range("M136").Value = -4656.8181
$ws.Range("N136").Value = -14017.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2448.2068
$ws.Range("I58").Value = 1887.9445
$ws.Range("J58").Value = 3365
$ws.Range("K58").Value = 1887.9445
$ws.Range("L58").Value = 3365
$ws.Range("M58").Value = -1684.9445
$ws.Range("N58").Value = -3771
$ws.Range("H62").Value = 4767.727
$ws.Range("I62").Value = 3271.6667
$ws.Range("J62").Value = 11500
$ws.Range("K62").Value = 3271.6667
$ws.Range("L62").Value = 11500
$ws.Range("M62").Value = -2647.6667
$ws.Range("N62").Value = -12748
$ws.Range("H65").Value = 4767.727
$ws.Range("I65").Value = 3271.6667
$ws.Range("J65").Value = 11500
$ws.Range("K65").Value = 16358.3335
$ws.Range("L65").Value = 57500
$ws.Range("M65").Value = -13238.3335
$ws.Range("N65").Value = -63740
$ws.Range("H99").Value = 5718.913
$ws.Range("I99").Value = 1291.3158
$ws.Range("K99").Value = 1291.3158
$ws.Range("M99").Value = 206.6841999999999
$ws.Range("H126").Value = 5718.913
$ws.Range("I126").Value = 1291.3158
$ws.Range("K126").Value = 3873.9474
$ws.Range("M126").Value = -1403.9474
$ws.Range("H132").Value = 4620
$ws.Range("I132").Value = 4439.8
$ws.Range("J132").Value = 4748.7144
$ws.Range("K132").Value = 13319.4
$ws.Range("L132").Value = 14246.1432
$ws.Range("M132").Value = -10789.4
$ws.Range("N132").Value = -19306.1432
$ws.Range("H134").Value = 2352.9285
$ws.Range("I134").Value = 2416.087
$ws.Range("K134").Value = 7248.261
$ws.Range("M134").Value = -4713.261
$ws.Range("H136").Value = 2448.2068
$ws.Range("I136").Value = 1887.9445
$ws.Range("J136").Value = 3365
$ws.Range("K136").Value = 5663.833500000001
$ws.Range("L136").Value = 10095
$ws.Range("M136").Value = -3113.833500000001
$ws.Range("N136").Value = -15195

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 891.1852
$ws.Range("I113").Value = 1057.5264
$ws.Range("J113").Value = 496.125
$ws.Range("K113").Value = 3172.5792
$ws.Range("L113").Value = 1488.375
$ws.Range("M113").Value = -1002.5792
$ws.Range("N113").Value = -5828.375
$ws.Range("H115").Value = 2476.1428
$ws.Range("I115").Value = 667.3333
$ws.Range("J115").Value = 3199.6667
$ws.Range("K115").Value = 2001.9999
$ws.Range("L115").Value = 9599.000100000001
$ws.Range("M115").Value = -826.9999
$ws.Range("N115").Value = -11949.0001
$ws.Range("H122").Value = 1035.8096
$ws.Range("I122").Value = 513.3158
$ws.Range("J122").Value = 5999.5
$ws.Range("K122").Value = 4619.8422
$ws.Range("L122").Value = 53995.5
$ws.Range("M122").Value = -2169.8422
$ws.Range("N122").Value = -58895.5
$ws.Range("H131").Value = 833.75
$ws.Range("J131").Value = 857.97894
$ws.Range("L131").Value = 2573.93682
$ws.Range("N131").Value = -12653.93682

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5810.1514
$ws.Range("I70").Value = 4846.3687
$ws.Range("K70").Value = 4846.3687
$ws.Range("M70").Value = -4576.3687
$ws.Range("H73").Value = 5810.1514
$ws.Range("I73").Value = 4846.3687
$ws.Range("K73").Value = 4846.3687
$ws.Range("M73").Value = -3910.3687
$ws.Range("H122").Value = 2400.225
$ws.Range("I122").Value = 1843.1482
$ws.Range("J122").Value = 3557.2307
$ws.Range("K122").Value = 5529.444600000001
$ws.Range("L122").Value = 10671.6921
$ws.Range("M122").Value = -3079.444600000001
$ws.Range("N122").Value = -15571.6921
$ws.Range("H126").Value = 2077.6667
$ws.Range("I126").Value = 1851.1
$ws.Range("J126").Value = 2360.875
$ws.Range("K126").Value = 5553.299999999999
$ws.Range("L126").Value = 7082.625
$ws.Range("M126").Value = -3083.299999999999
$ws.Range("N126").Value = -12022.625
$ws.Range("H132").Value = 2155.795
$ws.Range("I132").Value = 1632.037
$ws.Range("J132").Value = 3334.25
$ws.Range("K132").Value = 4896.111
$ws.Range("L132").Value = 10002.75
$ws.Range("M132").Value = -2366.111
$ws.Range("N132").Value = -15062.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 28812.895
$ws.Range("I40").Value = 43470.418
$ws.Range("J40").Value = 3685.7144
$ws.Range("K40").Value = 43470.418
$ws.Range("L40").Value = 3685.7144
$ws.Range("M40").Value = -43334.418
$ws.Range("N40").Value = -3957.7144
$ws.Range("H108").Value = 31485
$ws.Range("J108").Value = 31485
$ws.Range("L108").Value = 31485
$ws.Range("N108").Value = -39165
$ws.Range("H127").Value = 55751.8
$ws.Range("J127").Value = 55751.8
$ws.Range("L127").Value = 55751.8
$ws.Range("N127").Value = -65671.8
$ws.Range("H133").Value = 57499.5
$ws.Range("J133").Value = 57499.5
$ws.Range("L133").Value = 57499.5
$ws.Range("N133").Value = -62559.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""
$ws.Range("H132").Value = 502265
$ws.Range("I132").Value = 835316.4399999999
$ws.Range("J132").Value = 2687.875
$ws.Range("K132").Value = 2505949.32
$ws.Range("L132").Value = 8063.625
$ws.Range("M132").Value = -2503419.32
$ws.Range("N132").Value = -13123.625
